## klocki.xlsx edit: add a 5x10 coordinate-label grid (two 5x5 blocks) in F13:J17
## and L13:P17, mark B15 as a center ("ÅšR") cell, resize rows 13 and 14 to match
## the sheet's standard row height, and move the selection to F16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells already carrying the border/fill combinations the new grid
# needs (so we reuse existing style entries instead of minting new ones).
$styleSrc = @{
  2  = "B2"
  3  = "C2"
  4  = "F2"
  5  = "B3"
  6  = "E3"
  7  = "F3"
  8  = "B6"
  9  = "C6"
  11 = "J3"
  12 = "F4"
  55 = "C3"
  56 = "X3"
  57 = "D6"
  58 = "F6"
}

$targets = @(
  @{Cell="F13"; Style=2;  Text="0 0"},
  @{Cell="G13"; Style=3;  Text="0 1"},
  @{Cell="H13"; Style=3;  Text="0 2"},
  @{Cell="I13"; Style=3;  Text="0 3"},
  @{Cell="J13"; Style=4;  Text="0 4"},
  @{Cell="L13"; Style=2;  Text="4 5"},
  @{Cell="M13"; Style=3;  Text="5 5"},
  @{Cell="N13"; Style=3;  Text="6 5"},
  @{Cell="O13"; Style=3;  Text="7 5"},
  @{Cell="P13"; Style=4;  Text="8 5"},
  @{Cell="F14"; Style=5;  Text="1 0"},
  @{Cell="G14"; Style=55; Text="1 1"},
  @{Cell="H14"; Style=55; Text="1 2"},
  @{Cell="I14"; Style=6;  Text="1 3"},
  @{Cell="J14"; Style=7;  Text="1 4"},
  @{Cell="L14"; Style=5;  Text="4 6"},
  @{Cell="M14"; Style=55; Text="5 6"},
  @{Cell="N14"; Style=55; Text="6 6"},
  @{Cell="O14"; Style=6;  Text="7 6"},
  @{Cell="P14"; Style=7;  Text="8 6"},
  @{Cell="F15"; Style=5;  Text="2 0"},
  @{Cell="G15"; Style=11; Text="2 1"},
  @{Cell="H15"; Style=11; Text="2 2"},
  @{Cell="I15"; Style=11; Text="2 3"},
  @{Cell="J15"; Style=12; Text="2 4"},
  @{Cell="L15"; Style=5;  Text="4 7"},
  @{Cell="M15"; Style=11; Text="5 7"},
  @{Cell="N15"; Style=11; Text="6 7"},
  @{Cell="O15"; Style=11; Text="7 7"},
  @{Cell="P15"; Style=12; Text="8 7"},
  @{Cell="F16"; Style=5;  Text="3 0"},
  @{Cell="G16"; Style=6;  Text="3 1"},
  @{Cell="H16"; Style=55; Text="3 2"},
  @{Cell="I16"; Style=55; Text="3 3"},
  @{Cell="J16"; Style=56; Text="3 4"},
  @{Cell="L16"; Style=5;  Text="4 8"},
  @{Cell="M16"; Style=6;  Text="5 8"},
  @{Cell="N16"; Style=55; Text="6 8"},
  @{Cell="O16"; Style=55; Text="8 8"},
  @{Cell="P16"; Style=56; Text="7 8"},
  @{Cell="F17"; Style=8;  Text="4 0"},
  @{Cell="G17"; Style=9;  Text="4 1"},
  @{Cell="H17"; Style=57; Text="4 2"},
  @{Cell="I17"; Style=57; Text="4 3"},
  @{Cell="J17"; Style=58; Text="4 4"},
  @{Cell="L17"; Style=8;  Text="4 9"},
  @{Cell="M17"; Style=9;  Text="5 9"},
  @{Cell="N17"; Style=57; Text="6 9"},
  @{Cell="O17"; Style=57; Text="7 9"},
  @{Cell="P17"; Style=58; Text="8 9"}
)

foreach ($t in $targets) {
    $src = $ws.Range($styleSrc[[string]$t.Style])
    $dst = $ws.Range($t.Cell)
    $src.Copy()
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $dst.Value = $t.Text
}

# New rows 13/14 need the same row height as the rest of the sheet (created
# implicitly above without it).
$ws.Rows.Item(13).RowHeight = 27
$ws.Rows.Item(14).RowHeight = 27

# B15 becomes the grid's centre marker, same as the other tetromino grids.
$ws.Range("B15").Value = "ŚR"

$excel.CutCopyMode = $false

# Match the author's final selection / scroll position.
$ws.Range("F16").Select()
